$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$m.HeadersFooters.DateAndTime.Value = "2007-04-01"
